# Rename the "EDT Level 2 ... Functional Condition" crosswalk labels in
# column B to the new "EDT_Level 2 ..." / "EDT_..." naming scheme used by
# the updated Rank script, and populate the previously-blank "Side Channel"
# row with its own crosswalk name.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B5").Value  = "EDT_Level 2 Bed scour"
$ws.Range("B6").Value  = "EDT_Level 2 Benthic Richness"
$ws.Range("B8").Value  = "EDT_Level 2 Confinement: Artificial"
$ws.Range("B11").Value = "EDT_Level 2 Embeddedness"
$ws.Range("B12").Value = "EDT_Level 2 Fine Sediment"
$ws.Range("B18").Value = "EDT_Level 2 Flow: Inter-Annual Low Flow Variation"
$ws.Range("B19").Value = "EDT_Level 2 Flow: Inter-Annual High Flow Variation"
$ws.Range("B31").Value = "EDT_Level 2 Predation Risk"
$ws.Range("B32").Value = "EDT_Level 2 Riparian/stream interface"
$ws.Range("B38").Value = "EDT_Temperature: Daily Maximum"
$ws.Range("B40").Value = "EDT_Temperature: Food Effect"
$ws.Range("B43").Value = "EDT_Width"
$ws.Range("B44").Value = "EDT_Woody Debris"
$ws.Range("B36").Value = "EDT Level 2 Side Channel Functional Condition"

$ws.Range("B16").Select()
